$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.839.57"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "2.091.88"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.81"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.68"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0838"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "2.401.64"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.98"
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.796"
$ws.Range("E15").Value = "  +4.31%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "2.104.43"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "38.730.05"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.85"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.20"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.21"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("E28").Value = "  +8.74%  "
$ws.Range("E29").Value = "  +12.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.18"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  +6.10%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.43"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.20"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "1.541.75"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.81"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0225"
$ws.Range("E43").Value = "  +4.38%  "
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.66"
$ws.Range("E46").Value = "  +9.14%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "2.288.04"
$ws.Range("E51").Value = "  +2.46%  "
